$wb = $excel.ActiveWorkbook

# Column F ("想去人数") values were updated on both the "展览" sheet and the
# "全部类型" sheet (they contain duplicate data).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 2388
    $ws.Range("F4").Value = 439
    $ws.Range("F6").Value = 6493
    $ws.Range("F7").Value = 342
    $ws.Range("F8").Value = 128
}
